# Auto-generated Excel COM-interop edit script
# Applies the "cryptos list" update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text updates (value is not numeric-looking, safe to assign directly) ---
$ws.Range("D2").Value = "56.771.18"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "3.068.77"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("E6").Value = "  -3.52%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.067.48"
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("E9").Value = "  +1.95%  "
$ws.Range("E10").Value = "  +2.48%  "
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("D13").Value = "3.600.45"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("E16").Value = "  -2.13%  "
$ws.Range("D17").Value = "56.882.18"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "3.073.76"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("E19").Value = "  -3.73%  "
$ws.Range("E20").Value = "  -2.05%  "
$ws.Range("E21").Value = "  -1.76%  "
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("E26").Value = "  -2.89%  "
$ws.Range("E27").Value = "  -1.30%  "
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").Value = "0.0₃0853"
$ws.Range("E29").Value = "  -6.41%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("E33").Value = "  -7.90%  "
$ws.Range("E34").Value = "  -0.84%  "
$ws.Range("E35").Value = "  +6.08%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("E37").Value = "  -3.83%  "
$ws.Range("E38").Value = "  -3.00%  "
$ws.Range("E39").Value = "  -1.36%  "
$ws.Range("E40").Value = "  -1.76%  "
$ws.Range("E41").Value = "  -2.12%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("D45").Value = "2.376.47"
$ws.Range("E45").Value = "  +4.50%  "
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").Value = "3.110.36"
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("E50").Value = "  -3.76%  "
$ws.Range("E51").Value = "  -2.40%  "

# --- Numeric-looking text updates: force the cell to Text format first so Excel
#     keeps the literal string (e.g. "1.00", "0.999") instead of coercing it to a
#     number and losing the formatted trailing zeros / precision. ---
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "518.65"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "135.58"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.449"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "7.31"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.395"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "25.20"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0000160"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "5.84"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.39"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "7.80"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "345.76"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "68.10"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.495"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.166"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.24"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.86"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.85"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "20.70"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "4.88"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "158.92"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.13"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "5.96"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "25.59"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0649"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.57"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "3.99"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.686"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "36.54"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.952"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "5.92"
